$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 63

# Text-like columns (Date, Time, Weekday, Week) must stay as plain text,
# matching the inlineStr storage used throughout the sheet - force text
# number format while assigning, then restore the default "Normal" style
# so no explicit style index lingers on the new cells.
function Set-TextCell($r, $c, $val) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $row 1 "2025-02-11"
Set-TextCell $row 2 "08:53:26"
Set-TextCell $row 3 "Tuesday"
Set-TextCell $row 4 "06"

$ws.Cells.Item($row, 5).Value = 127563
$ws.Cells.Item($row, 6).Value = 141958
$ws.Cells.Item($row, 7).Value = 169016
$ws.Cells.Item($row, 8).Value = 158417
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 144253
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191761
$ws.Cells.Item($row, 14).Value = 115057
$ws.Cells.Item($row, 15).Value = 44855
$ws.Cells.Item($row, 16).Value = 28487
$ws.Cells.Item($row, 17).Value = 64590
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42488
$ws.Cells.Item($row, 20).Value = -1
